$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2023976666666667
$ws.Range("N2").Value = 0.6071930000000001
$ws.Range("O2").Value = 0.03663970451354832
$ws.Range("P2").Value = 0.03663970451354832
$ws.Range("Q2").Value = 1.566766207199
$ws.Range("R2").Value = 14.100895864791
$ws.Range("S2").Value = 0.0180651584362809
$ws.Range("T2").Value = 0.0180651584362809

# Row 3
$ws.Range("G3").Value = 7.741029
$ws.Range("H3").Value = 23.223087
$ws.Range("I3").Value = 0.4930486933812723
$ws.Range("J3").Value = 0.4930486933812723
$ws.Range("O3").Value = 0.08641717548188978
$ws.Range("P3").Value = 0.08641717548188979
$ws.Range("Q3").Value = 3.695322112015
$ws.Range("R3").Value = 33.257899008135
$ws.Range("S3").Value = 0.04260787545704587
$ws.Range("T3").Value = 0.04260787545704588

# Row 4
$ws.Range("G4").Value = 7.741029
$ws.Range("H4").Value = 23.223087
$ws.Range("I4").Value = 0.4930486933812723
$ws.Range("J4").Value = 0.4930486933812723
$ws.Range("O4").Value = 0.876943120004562
$ws.Range("P4").Value = 0.876943120004562
$ws.Range("Q4").Value = 37.499342975071
$ws.Range("R4").Value = 337.4940867756389
$ws.Range("S4").Value = 0.4323756594879455
$ws.Range("T4").Value = 0.4323756594879455

# Row 5
$ws.Range("I5").Value = 0.0194007766416684
$ws.Range("J5").Value = 0.0194007766416684
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2023976666666667
$ws.Range("N5").Value = 0.6071930000000001
$ws.Range("O5").Value = 0.03663970451354832
$ws.Range("P5").Value = 0.03663970451354832
$ws.Range("Q5").Value = 0.06165005940311112
$ws.Range("R5").Value = 0.5548505346280002
$ws.Range("S5").Value = 0.0007108387234840805
$ws.Range("T5").Value = 0.0007108387234840805

# Row 6
$ws.Range("I6").Value = 0.0194007766416684
$ws.Range("J6").Value = 0.0194007766416684
$ws.Range("O6").Value = 0.08641717548188978
$ws.Range("P6").Value = 0.08641717548188979
$ws.Range("S6").Value = 0.001676560319528006
$ws.Range("T6").Value = 0.001676560319528007

# Row 7
$ws.Range("I7").Value = 0.0194007766416684
$ws.Range("J7").Value = 0.0194007766416684
$ws.Range("O7").Value = 0.876943120004562
$ws.Range("P7").Value = 0.876943120004562
$ws.Range("S7").Value = 0.01701337759865632
$ws.Range("T7").Value = 0.01701337759865632

# Row 8
$ws.Range("G8").Value = 7.654706000000001
$ws.Range("I8").Value = 0.4875505299770593
$ws.Range("J8").Value = 0.4875505299770593
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2023976666666667
$ws.Range("N8").Value = 0.6071930000000001
$ws.Range("O8").Value = 0.03663970451354832
$ws.Range("P8").Value = 0.03663970451354832
$ws.Range("Q8").Value = 1.549294633419334
$ws.Range("R8").Value = 13.943651700774
$ws.Range("S8").Value = 0.01786370735378334
$ws.Range("T8").Value = 0.01786370735378334

# Row 9
$ws.Range("G9").Value = 7.654706000000001
$ws.Range("I9").Value = 0.4875505299770593
$ws.Range("J9").Value = 0.4875505299770593
$ws.Range("O9").Value = 0.08641717548188978
$ws.Range("P9").Value = 0.08641717548188979
$ws.Range("R9").Value = 32.88702820839001
$ws.Range("S9").Value = 0.0421327397053159
$ws.Range("T9").Value = 0.04213273970531591

# Row 10
$ws.Range("G10").Value = 7.654706000000001
$ws.Range("I10").Value = 0.4875505299770593
$ws.Range("J10").Value = 0.4875505299770593
$ws.Range("O10").Value = 0.876943120004562
$ws.Range("P10").Value = 0.876943120004562
$ws.Range("R10").Value = 333.730568766246
$ws.Range("S10").Value = 0.4275540829179602
$ws.Range("T10").Value = 0.4275540829179602
